# swot_analyisis.xlsx – finish the SWOT cross-analysis table and drop a
# stray test entry from the Strengths/Opportunities lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove leftover "test" rows from the two list columns -----------------
$ws.Range("D4").Value = "Zu wenig Kundschaft (da abgelegen)"
$ws.Range("C5").Value = "Neugierige Kunden"
$ws.Range("D6").Value = "zu teure Miete"

# rename a strength entry
$ws.Range("B11").Value = "Qualitativ gute Musiker (Jungmusiker)"

# --- fill in the SWOT cross-analysis (Strengths / Weaknesses vs. Opportunities / Threats) --
# header row for the "Strengths" sub-table
$ws.Range("D8").Value = "Welche Stärken treffen auf welche Gefahren?"
$ws.Range("D8").Font.Bold = $true

$ws.Range("C9").Value = "Durch Treue Kunden gewinnen wir die Unterstützung des Dorfes"
$ws.Range("D9").Value = "Durch gute Musiker + Preisleistungsverhältniss mindert sich der Konkurrentdruck"

$ws.Range("C10").Value = "Durch die Parkölätze sind auch grössere Evente möglich"
$ws.Range("D10").Value = "Durch treue Kunden können wir die Kosten der Miete begleichen"

$ws.Range("C11").Value = "Durch gute Musiker kommen neugiereige Kunden "
$ws.Range("D11").Value = "Durch ausgewählte Jungmusiker haben wir geringere Kosten"

$ws.Range("C12").Value = "Durch das gute Preisleistungsverhältnis bekommt man Stammkunden"

# header row for the "Weaknesses" sub-table
$ws.Range("C14").Value = "Welche Schwächen treffen auf welche Möglichkeiten?"
$ws.Range("C14").Font.Bold = $true
$ws.Range("D14").Value = "Welche Schwächen treffen auf welche Gefahren?"
$ws.Range("D14").Font.Bold = $true

$ws.Range("C15").Value = "Die unglücklichen Nachbaren von der Idee überzeugen -> unterstützen uns"
$ws.Range("D15").Value = "evt. Umzug in umsatzstärkeres Gebiet"

$ws.Range("C16").Value = "Abgelegener Ort ist ideal für grössere Veranstaltungen, bringt dem Dorf Gewinne ein"
$ws.Range("D16").Value = "Durch genügend Parkplätze wird die Abgelegenheit kompensiert"

# leave the current selection where the author left off editing
$ws.Range("C18").Select() | Out-Null
